# "New photos since June 2022 until now, some homepage reorganization"
#
# The sheet is a small helper table used to hand-build masonry <div><img>
# markup: column A holds the literal (syntax-highlighted) HTML prefix,
# column B holds a bare filename pasted in by hand, column C holds the
# HTML suffix (alt text + closing tags, also syntax-highlighted), and
# column D concatenates A/B/C into the final markup string.
#
# This edit swaps the old "Portraits" rows (2-18) for four new "Pretty
# Things" rows (2-5), blanks rows 6-18 back to the same empty "template"
# state as the rows below them, and extends the sheet's pre-built blank
# template rows from 210 down to 260.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helpers for writing "rich text" cells (multiple colored runs inside a
# single cell) the same way the syntax-highlighted A/C columns use.
# ---------------------------------------------------------------------
function ConvertTo-ComColor {
    param([string]$hex)
    if ($hex.Length -gt 6) { $hex = $hex.Substring($hex.Length - 6) }
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# $runs is an array of 2-item arrays: @(<hex color or $null>, <text>)
# A $null color leaves that run on the cell's default font color.
function Set-RichText {
    param($cell, $runs)
    $text = ""
    foreach ($run in $runs) { $text += $run[1] }
    $cell.Value = $text

    $pos = 1
    foreach ($run in $runs) {
        $hex = $run[0]
        $len = $run[1].Length
        if ($hex) {
            $chars = $cell.Characters($pos, $len)
            $chars.Font.Color = (ConvertTo-ComColor $hex)
        }
        $pos += $len
    }
}

# ---------------------------------------------------------------------
# New data rows (2-5): "Pretty Things" folder
# ---------------------------------------------------------------------

# Column A prefix is identical on every data row.
$runsA = @(
    @($null, '<'),
    @("FFEF596F", "div"),
    @("FFBBBBBB", " "),
    @("FFD19A66", "class"),
    @("FFBBBBBB", "="),
    @("FF89CA78", '"masonryImage"'),
    @("FFBBBBBB", "> <"),
    @("FFFF0000", "img"),
    @("FFBBBBBB", " "),
    @("FFBF8F00", "src="),
    @("FFA9D08E", '"images'),
    @("FFBBBBBB", "/Pretty Things/")
)

Set-RichText $ws.Range("A2") $runsA
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4104)
$excel.CutCopyMode = 0

function Get-AltRuns {
    param([string]$altText)
    return @(
        @($null, '" '),
        @("FFBF8F00", "alt="),
        @("FFA9D08E", ('"' + $altText + '"')),
        @("FF89CA78", " "),
        @("FFD0CECE", "/></"),
        @("FFFF0000", "div>")
    )
}

Set-RichText $ws.Range("C2") (Get-AltRuns "Fall colors in Great Barrington Massachussets in Fall of 2022.")
Set-RichText $ws.Range("C3") (Get-AltRuns "Fall colors on a winding road in Great Barrington Massachussets in Fall of 2022.")
Set-RichText $ws.Range("C4") (Get-AltRuns "Photos from California during Winter 2022 and Summer 2023.")

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("B2").Value = "DSCF0846-Editweb.jpg"
$ws.Range("B3").Value = "DSCF0851web.jpg"
$ws.Range("B4").Value = "DSCF1302web.jpg"
$ws.Range("B5").Value = "DSCF3793web.jpg"

# D2:D5 rebuild the concatenated markup from A/B/C.
$ws.Range("D2").Formula = "=CONCATENATE(A2,B2,C2)"
$ws.Range("D3").Formula = "=CONCATENATE(A3,B3,C3)"
$ws.Range("D4").Formula = "=CONCATENATE(A4,B4,C4)"
$ws.Range("D5").Formula = "=CONCATENATE(A5,B5,C5)"

# Row 3's wrapped alt text is one line shorter than rows 2/4/5.
$ws.Rows("3").RowHeight = 48

# ---------------------------------------------------------------------
# Rows 6-18 revert to the empty "template" state (they used to hold the
# tail end of the old "Portraits" list).
# ---------------------------------------------------------------------
$ws.Range("A6:D18").ClearContents()
$ws.Rows("6:18").AutoFit()

# ---------------------------------------------------------------------
# Grow the blank template rows at the bottom from 210 to 260. Row 20 is
# an already-blank template row (style only, no B cell) -- copy its
# per-column formatting down into the freshly added rows.
# ---------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("A211:A260").PasteSpecial(-4122)
$ws.Range("C20").Copy()
$ws.Range("C211:C260").PasteSpecial(-4122)
$ws.Range("D20").Copy()
$ws.Range("D211:D260").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# View state: selection moved to C3, no frozen/scrolled top-left cell.
# ---------------------------------------------------------------------
$ws.Range("C3").Select()
